$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the UserDashboardAccountId column (column J) ---
$ws.Columns("J:J").Delete()

# --- Update header row (unchanged text, columns A:I stay the same) ---

# --- Row 2: UserAccountId=1 ---
$ws.Range("B2").Value = 1066898235
$ws.Range("C2").Value = "Bank"
$ws.Range("D2").Value = "Test User Welcome"
$ws.Range("E2").Value = "Savings Account "
$ws.Range("F2").Value = "SA"
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 10
$ws.Range("I2").Value = 1

# --- Row 3 ---
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "1066898235"
$ws.Range("C3").Value = "Bank"
$ws.Range("D3").Value = "Tom Savings"
$ws.Range("E3").Value = "Savings Account "
$ws.Range("F3").Value = "SA"
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 10
$ws.Range("I3").Value = 11

# --- Row 4 ---
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "1602334536"
$ws.Range("C4").Value = "Bank"
$ws.Range("D4").Value = "Beth Savings"
$ws.Range("E4").Value = "Savings Account "
$ws.Range("F4").Value = "SA"
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 10
$ws.Range("I4").Value = 12

# --- Row 5 ---
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "1633281132"
$ws.Range("C5").Value = "Bank"
$ws.Range("D5").Value = "Danny Savings"
$ws.Range("E5").Value = "Savings Account "
$ws.Range("F5").Value = "SA"
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 10
$ws.Range("I5").Value = 13

# --- Row 6 ---
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "4922130000784779"
$ws.Range("C6").Value = "Card"
$ws.Range("D6").Value = "Test User Card"
$ws.Range("E6").Value = "Savings Account "
$ws.Range("F6").Value = "CC"
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 10
$ws.Range("I6").Value = 1

# --- Row 7 ---
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "5898460761682640"
$ws.Range("C7").Value = "Card"
$ws.Range("D7").Value = "Beth Card"
$ws.Range("E7").Value = "Savings Account "
$ws.Range("F7").Value = "CC"
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 10
$ws.Range("I7").Value = 11

# --- Extra blank rows 9:12 (row 8 already blank/styled) ---
$ws.Range("A9").Value = $null
$ws.Range("A10").Value = $null
$ws.Range("A11").Value = $null
$ws.Range("A12").Value = $null

# --- Sheet view changes ---
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("H13").Select()
